$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 667.5
$ws.Range("I19").Value = 695
$ws.Range("J19").Value = 658.3333
$ws.Range("K19").Value = 695
$ws.Range("L19").Value = 658.3333
$ws.Range("M19").Value = -520
$ws.Range("N19").Value = -1008.3333
$ws.Range("H125").Value = 3678.75
$ws.Range("J125").Value = 3290.0715
$ws.Range("L125").Value = 29610.6435
$ws.Range("N125").Value = -34530.6435
$ws.Range("H127").Value = 767.2
$ws.Range("I127").Value = 634
$ws.Range("K127").Value = 1902
$ws.Range("M127").Value = 3058
$ws.Range("H129").Value = 1093.7755
$ws.Range("J129").Value = 1221.4147
$ws.Range("L129").Value = 3664.2441
$ws.Range("N129").Value = -13664.2441
$ws.Range("H137").Value = 2101.7856
$ws.Range("I137").Value = 1419.2941
$ws.Range("J137").Value = 3156.5454
$ws.Range("K137").Value = 4257.8823
$ws.Range("L137").Value = 9469.636200000001
$ws.Range("M137").Value = -1707.8823
$ws.Range("N137").Value = -14569.6362
$ws.Range("H138").Value = 3022.5806
$ws.Range("I138").Value = 1637.1538
$ws.Range("J138").Value = 4023.1667
$ws.Range("K138").Value = 4911.4614
$ws.Range("L138").Value = 12069.5001
$ws.Range("M138").Value = 228.5385999999999
$ws.Range("N138").Value = -22349.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 500
$ws.Range("I45").Value = 500
$ws.Range("K45").Value = 500
$ws.Range("M45").Value = -123
$ws.Range("H57").Value = 3000
$ws.Range("I57").Value = 3000
$ws.Range("K57").Value = 3000
$ws.Range("M57").Value = -2516
$ws.Range("H61").Value = 2695.762
$ws.Range("I61").Value = 2140.7334
$ws.Range("J61").Value = 4083.3333
$ws.Range("K61").Value = 2140.7334
$ws.Range("L61").Value = 4083.3333
$ws.Range("M61").Value = -1928.7334
$ws.Range("N61").Value = -4507.3333
$ws.Range("H97").Value = 551.7646999999999
$ws.Range("I97").Value = 561.6667
$ws.Range("J97").Value = 528
$ws.Range("K97").Value = 561.6667
$ws.Range("L97").Value = 528
$ws.Range("M97").Value = -65.66669999999999
$ws.Range("N97").Value = -1520
$ws.Range("H107").Value = 40000
$ws.Range("J107").Value = 40000
$ws.Range("L107").Value = 40000
$ws.Range("N107").Value = -47680
$ws.Range("H109").Value = 34000
$ws.Range("J109").Value = 34000
$ws.Range("L109").Value = 34000
$ws.Range("N109").Value = -36774
$ws.Range("H122").Value = 2915.4546
$ws.Range("I122").Value = 3480.75
$ws.Range("J122").Value = 1408
$ws.Range("K122").Value = 10442.25
$ws.Range("L122").Value = 4224
$ws.Range("M122").Value = -7992.25
$ws.Range("N122").Value = -9124
$ws.Range("H123").Value = 26415.2
$ws.Range("J123").Value = 26415.2
$ws.Range("L123").Value = 26415.2
$ws.Range("N123").Value = -36215.2
$ws.Range("H136").Value = 2695.762
$ws.Range("I136").Value = 2140.7334
$ws.Range("J136").Value = 4083.3333
$ws.Range("K136").Value = 6422.2002
$ws.Range("L136").Value = 12249.9999
$ws.Range("M136").Value = -3872.2002
$ws.Range("N136").Value = -17349.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 92353.82000000001
$ws.Range("I20").Value = 144384.72
$ws.Range("J20").Value = 1299.75
$ws.Range("K20").Value = 144384.72
$ws.Range("L20").Value = 1299.75
$ws.Range("M20").Value = -144137.72
$ws.Range("N20").Value = -1793.75
$ws.Range("H22").Value = 5195.2
$ws.Range("I22").Value = 5195.2
$ws.Range("K22").Value = 5195.2
$ws.Range("M22").Value = -5022.2
$ws.Range("H94").Value = 1112.909
$ws.Range("I94").Value = 910.82355
$ws.Range("J94").Value = 1800
$ws.Range("K94").Value = 910.82355
$ws.Range("L94").Value = 1800
$ws.Range("M94").Value = -459.82355
$ws.Range("N94").Value = -2702
$ws.Range("H99").Value = 4200
$ws.Range("I99").Value = 8000
$ws.Range("J99").Value = 3250
$ws.Range("K99").Value = 8000
$ws.Range("L99").Value = 3250
$ws.Range("M99").Value = -6502
$ws.Range("N99").Value = -6246

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2483.8262
$ws.Range("I31").Value = 1302
$ws.Range("J31").Value = 5483.846
$ws.Range("K31").Value = 1302
$ws.Range("L31").Value = 5483.846
$ws.Range("M31").Value = -1007
$ws.Range("N31").Value = -6073.846
$ws.Range("H34").Value = 2483.8262
$ws.Range("I34").Value = 1302
$ws.Range("J34").Value = 5483.846
$ws.Range("K34").Value = 1302
$ws.Range("L34").Value = 5483.846
$ws.Range("M34").Value = -1100
$ws.Range("N34").Value = -5887.846
$ws.Range("H58").Value = 1599.6129
$ws.Range("I58").Value = 1227.8636
$ws.Range("J58").Value = 2508.3333
$ws.Range("K58").Value = 1227.8636
$ws.Range("L58").Value = 2508.3333
$ws.Range("M58").Value = -1024.8636
$ws.Range("N58").Value = -2914.3333
$ws.Range("H122").Value = 2906.4
$ws.Range("I122").Value = 2105.2856
$ws.Range("J122").Value = 4775.6665
$ws.Range("K122").Value = 6315.8568
$ws.Range("L122").Value = 14326.9995
$ws.Range("M122").Value = -3865.8568
$ws.Range("N122").Value = -19226.9995
$ws.Range("H134").Value = 1922.6216
$ws.Range("I134").Value = 1528.92
$ws.Range("J134").Value = 2742.8333
$ws.Range("K134").Value = 4586.76
$ws.Range("L134").Value = 8228.499899999999
$ws.Range("M134").Value = -2051.76
$ws.Range("N134").Value = -13298.4999
$ws.Range("H136").Value = 1599.6129
$ws.Range("I136").Value = 1227.8636
$ws.Range("J136").Value = 2508.3333
$ws.Range("K136").Value = 3683.5908
$ws.Range("L136").Value = 7524.999899999999
$ws.Range("M136").Value = -1133.5908
$ws.Range("N136").Value = -12624.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6083.9473
$ws.Range("J70").Value = 6199.2856
$ws.Range("L70").Value = 6199.2856
$ws.Range("N70").Value = -6739.2856
$ws.Range("H73").Value = 6083.9473
$ws.Range("J73").Value = 6199.2856
$ws.Range("L73").Value = 6199.2856
$ws.Range("N73").Value = -8071.2856
$ws.Range("H102").Value = 3002.6428
$ws.Range("I102").Value = 3192.1428
$ws.Range("J102").Value = 2813.1428
$ws.Range("K102").Value = 3192.1428
$ws.Range("L102").Value = 2813.1428
$ws.Range("M102").Value = -1570.1428
$ws.Range("N102").Value = -6057.1428
$ws.Range("H123").Value = 18300.143
$ws.Range("J123").Value = 18300.143
$ws.Range("L123").Value = 18300.143
$ws.Range("N123").Value = -23200.143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 9000
$ws.Range("I30").Value = 11000
$ws.Range("J30").Value = 5000
$ws.Range("K30").Value = 11000
$ws.Range("L30").Value = 5000
$ws.Range("M30").Value = -10892
$ws.Range("N30").Value = -5216
$ws.Range("H40").Value = 3375.75
$ws.Range("I40").Value = 3600
$ws.Range("J40").Value = 2703
$ws.Range("K40").Value = 3600
$ws.Range("L40").Value = 2703
$ws.Range("M40").Value = -3464
$ws.Range("N40").Value = -2975
$ws.Range("H82").Value = 1950.7693
$ws.Range("I82").Value = 1419.4546
$ws.Range("J82").Value = 2340.4
$ws.Range("K82").Value = 1419.4546
$ws.Range("L82").Value = 2340.4
$ws.Range("M82").Value = -1058.4546
$ws.Range("N82").Value = -3062.4
$ws.Range("H85").Value = 1950.7693
$ws.Range("I85").Value = 1419.4546
$ws.Range("J85").Value = 2340.4
$ws.Range("K85").Value = 1419.4546
$ws.Range("L85").Value = 2340.4
$ws.Range("M85").Value = -171.4546
$ws.Range("N85").Value = -4836.4
$ws.Range("H104").Value = 20567.5
$ws.Range("J104").Value = 20567.5
$ws.Range("L104").Value = 20567.5
$ws.Range("N104").Value = -27555.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 64456.367
$ws.Range("J81").Value = 7943.8184
$ws.Range("L81").Value = 15887.6368
$ws.Range("N81").Value = -18009.6368
$ws.Range("H84").Value = 64456.367
$ws.Range("J84").Value = 7943.8184
$ws.Range("L84").Value = 79438.18400000001
$ws.Range("N84").Value = -90046.18400000001
$ws.Range("H126").Value = 4622.25
$ws.Range("I126").Value = 4868.636
$ws.Range("J126").Value = 4080.2
$ws.Range("K126").Value = 14605.908
$ws.Range("L126").Value = 12240.6
$ws.Range("M126").Value = -12135.908
$ws.Range("N126").Value = -17180.6
$ws.Range("H131").Value = 59999.5
$ws.Range("J131").Value = 59999.5
$ws.Range("L131").Value = 59999.5
$ws.Range("N131").Value = -70079.5
$ws.Range("H132").Value = 3662.963
$ws.Range("I132").Value = 3675.25
$ws.Range("J132").Value = 3657.7896
$ws.Range("K132").Value = 11025.75
$ws.Range("L132").Value = 10973.3688
$ws.Range("M132").Value = -8495.75
$ws.Range("N132").Value = -16033.3688
